# Apply updates to the "想去人数" (F) and "最低票价" (G) columns
# on the "展览" and "全部类型" worksheets, per the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value = 133
$ws1.Range("F4").Value = 2092
$ws1.Range("G4").Value = 55.2
$ws1.Range("F5").Value = 369
$ws1.Range("F6").Value = 646
$ws1.Range("F8").Value = 2079
$ws1.Range("F9").Value = 10746
$ws1.Range("F15").Value = 7595
$ws1.Range("F17").Value = 724
$ws1.Range("F18").Value = 698
$ws1.Range("F20").Value = 3346

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value = 133
$ws4.Range("F4").Value = 2092
$ws4.Range("G4").Value = 55.2
$ws4.Range("F5").Value = 369
$ws4.Range("F6").Value = 646
$ws4.Range("F9").Value = 2079
$ws4.Range("F12").Value = 10746
$ws4.Range("F18").Value = 7596
$ws4.Range("F20").Value = 724
$ws4.Range("F21").Value = 702
$ws4.Range("F23").Value = 3346
